# daily-log as pandas was not working to update excel, change to list of rows
# Replace the pandas-driven daily log update with an explicit list-of-rows
# write into the "May 01" sheet: clear the stray blank separator cells and
# append the day's food entries as plain rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("May 01")

# The old pandas-based update left behind empty placeholder cells on the
# blank separator rows (2 and 6) - clear them out now that rows are
# written directly.
$ws.Range("A2").ClearContents()
$ws.Range("A6").ClearContents()

# Food log entries for the day, as a plain list of rows instead of a
# pandas DataFrame: Food, Amount, Unit, Calories, Protein, Carbs, Fats
$foodRows = @(
    @("beer", 500, "ml", 215, 2.3, 18, 0),
    @("wine", 1, "unit", 162, 0.14, 4, 0),
    @("protein bar", 1, "unit", 210, 20, 26, 7)
)

$startRow = 8
for ($i = 0; $i -lt $foodRows.Count; $i++) {
    $rowValues = $foodRows[$i]
    $targetRow = $startRow + $i
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = $rowValues[$col - 1]
    }
}

# Leave the sheet active with the cell just past the new data selected,
# matching where the user's cursor ended up after the update.
$ws.Activate() | Out-Null
$ws.Range("J9").Select() | Out-Null
